$wb = $excel.ActiveWorkbook

# --- Fix the PIN MAPPING sheet frozen-pane scroll position (A34 -> A2), keep selection at D64 ---
$pm = $wb.Worksheets.Item("PIN MAPPING")
$pm.Activate()
$pm.Range("A2").Select()
$pm.Range("D64").Select()

# --- CHIPS sheet: populate manufacturer / part-number / description / datasheet-link columns ---
$ws = $wb.Worksheets.Item("CHIPS")
$ws.Activate()

# Row 16 used to carry its part number in column B; move it to column C
$ws.Range("B16").ClearContents()

$ws.Range("B1").Value = "ST"
$ws.Range("C1").Value = "uln2803a"
$ws.Range("D1").Value = "DarlingtonTransistor Arrays"
$ws.Range("E1").Value = "https://www.ti.com/lit/ds/symlink/uln2803a.pdf?ts=1645366810061&ref_url=https%253A%252F%252Fwww.google.com%252F"
$ws.Range("B2").Value = "TI"
$ws.Range("C2").Value = "sn74hc259n"
$ws.Range("B3").Value = "TI"
$ws.Range("C3").Value = "sn74hc02n"
$ws.Range("B4").Value = "Motorola"
$ws.Range("C4").Value = "LM324N"
$ws.Range("E4").Value = "https://pdf1.alldatasheet.com/datasheet-pdf/view/22756/STMICROELECTRONICS/LM324N.html"
$ws.Range("B5").Value = "TI"
$ws.Range("C5").Value = "ADC0808N"
$ws.Range("B6").Value = "Ti"
$ws.Range("C6").Value = "SN74HC74N"
$ws.Range("B7").Value = "TI"
$ws.Range("C7").Value = "SN74HC02N"
$ws.Range("B8").Value = "TI"
$ws.Range("C8").Value = "SN74HC74N"
$ws.Range("B9").Value = "Motorola"
$ws.Range("C9").Value = "mc14516bcp"
$ws.Range("B10").Value = "Motorola"
$ws.Range("C10").Value = "mc14516bcp"
$ws.Range("B11").Value = "Motorola"
$ws.Range("C11").Value = "MC74HC7266N"
$ws.Range("B12").Value = "Motorola"
$ws.Range("C12").Value = "MC74HC7266N"
$ws.Range("B13").Value = "Motorola"
$ws.Range("C13").Value = "MC74HC14"
$ws.Range("B14").Value = "TI"
$ws.Range("C14").Value = "SN74HC20N"
$ws.Range("B15").Value = "TI"
$ws.Range("C15").Value = "SN74HC138N"
$ws.Range("C16").Value = "74hc139"
$ws.Range("B17").Value = "Motorola"
$ws.Range("C17").Value = "MC14516BCP"
$ws.Range("B18").Value = "Motorola"
$ws.Range("C18").Value = "MC14516BCP"
$ws.Range("C19").Value = "MM14516BCN"
$ws.Range("B20").Value = "TI"
$ws.Range("C20").Value = "SN74HC573N"
$ws.Range("B21").Value = "TI"
$ws.Range("C21").Value = "SN74HC573N"
$ws.Range("B22").Value = "Motorola"
$ws.Range("C22").Value = "74HC573"
$ws.Range("B23").Value = "Motorola"
$ws.Range("C23").Value = "MC14516BCP"
$ws.Range("B24").Value = "Motorola"
$ws.Range("C24").Value = "MC14516BCP"
$ws.Range("B25").Value = "TI"
$ws.Range("C25").Value = "SN74HC74N"
$ws.Range("B26").Value = "TI"
$ws.Range("C26").Value = "SN74HC74N"
$ws.Range("B27").Value = "TI"
$ws.Range("C27").Value = "SN74HC74N"
$ws.Range("B28").Value = "Motorola"
$ws.Range("C28").Value = "MC74HC08"
$ws.Range("B29").Value = "TI"
$ws.Range("C29").Value = "SN74HC03N"
$ws.Range("B30").Value = "TI"
$ws.Range("C30").Value = "SN74HC74N"
$ws.Range("B31").Value = "TI"
$ws.Range("C31").Value = "SN74HC04N"
$ws.Range("B32").Value = "Toshiba"
$ws.Range("C32").Value = "TMPZ84C00AP"
$ws.Range("B33").Value = "Motorola"
$ws.Range("C33").Value = "MCM6064P12"
$ws.Range("B34").Value = "TI"
$ws.Range("C34").Value = "TMS 27C64-IJL"
$ws.Range("B35").Value = "TI"
$ws.Range("C35").Value = "TMS 27C64-IJL"
$ws.Range("B36").Value = "Zilog"
$ws.Range("C36").Value = "Z84C3006PEC"
$ws.Range("D36").Value = "Z80 CTC"
$ws.Range("B37").Value = "Zilog"
$ws.Range("C37").Value = "Z84C4004PEC"
$ws.Range("D37").Value = "Z80 SI0/0"
$ws.Range("B38").Value = "TI"
$ws.Range("C38").Value = "SN74HC157N"
$ws.Range("B39").Value = "TI"
$ws.Range("C39").Value = "SN74HC4060N"
$ws.Range("B40").Value = "Motorola"
$ws.Range("C40").Value = "MC1489P"
$ws.Range("B41").Value = "Motorola"
$ws.Range("C41").Value = "MC1488P"
$ws.Range("B42").Value = "Motorola"
$ws.Range("C42").Value = "MC1489P"

# Hyperlink the datasheet URL entered in column E1 (E4 is left as plain text)
$ws.Hyperlinks.Add($ws.Range("E1"), $ws.Range("E1").Value) | Out-Null

# Widen columns C and D to fit the new text
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 23

# Restore the view: scroll so row 18 is at the top, leave the cursor where editing left off
$ws.Range("A18").Select()
$ws.Range("B43").Select()
